$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (StudentID -> user_id, Score -> score)
$ws.Range("A1").Value = "user_id"
$ws.Range("B1").Value = "score"

# Update score values in column B (rows 4-23 changed per diff; others unchanged but set anyway for safety)
$scores = @{
    2 = 7
    3 = 8
    4 = 6.5
    5 = 8
    6 = 9
    7 = 10
    8 = 8
    9 = 6
    10 = 8.5
    11 = 7.5
    12 = 7
    13 = 9.5
    14 = 8.1999999999999993
    15 = 7.5
    16 = 6
    17 = 8
    18 = 9
    19 = 8.5
    20 = 7
    21 = 6
    22 = 5
    23 = 8.5
    24 = 9
    25 = 7
    26 = 9
    27 = 7
    28 = 10
}

foreach ($row in $scores.Keys) {
    $ws.Cells.Item($row, 2).Value = $scores[$row]
}

# Update the active selection to B1 as reflected in the saved worksheet view
$ws.Range("B1").Select()
